$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the style/formatting used by the other
# header cells (B1:G1) -- copy G1's formatting (font/border/alignment) onto H1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add new numeric data for the "Save" column, matching rows 2 and 3
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
